$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated save_data: column G ("K") now holds true strikeout counts
# (previously it held a "Strike#" pitch-count-derived value). The table
# below maps each data row (2-69) to its recalculated K value.
$kValues = @{
    2  = 0;  3  = 0;  4  = 1;  5  = 1;  6  = 0;  7  = 0;  8  = 0;  9  = 1;
    10 = 1;  11 = 4;  12 = 2;  13 = 0;  14 = 1;  15 = 1;  16 = 1;  17 = 2;
    18 = 1;  19 = 1;  20 = 2;  21 = 1;  22 = 1;  23 = 3;  24 = 1;  25 = 0;
    26 = 1;  27 = 0;  28 = 3;  29 = 0;  30 = 0;  31 = 0;  32 = 1;  33 = 0;
    34 = 3;  35 = 0;  36 = 2;  37 = 2;  38 = 0;  39 = 4;  40 = 1;  41 = 1;
    42 = 3;  43 = 1;  44 = 2;  45 = 0;  46 = 1;  47 = 2;  48 = 1;  49 = 1;
    50 = 2;  51 = 1;  52 = 1;  53 = 2;  54 = 0;  55 = 4;  56 = 0;  57 = 2;
    58 = 1;  59 = 0;  60 = 1;  61 = 3;  62 = 1;  63 = 0;  64 = 2;  65 = 2;
    66 = 0;  67 = 1;  68 = 2;  69 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
